$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Review date" column is column B; dates were auto-generated with the
# day-of-month fixed at "05" except for a handful of genuinely distinct
# dates. Re-generation bumped every "-05" day to "-17" while leaving any
# other day value untouched.
$lastRow = $ws.UsedRange.Rows.Count()

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $oldStr = [string]$cell.Value()

    if ($oldStr.EndsWith("-05")) {
        $newStr = $oldStr.Substring(0, $oldStr.Length - 2) + "17"

        # Force the cell to stay text (it already stores a text date like
        # "2011-03-05") instead of letting Excel reinterpret the new
        # string as a real date serial, then drop back to the default
        # style so no stray formatting is left behind.
        $cell.NumberFormat = "@"
        $cell.Value = $newStr
        $cell.Style = "Normal"
    }
}
